$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values per row (columns G,H,I,J = ligand stats; M,N,O,P = receptor stats; Q,R,S,T = edge weights)
$updates = @{
    "G2" = 25.3659715
    "H2" = 50.731943
    "I2" = 0.1464818225604603
    "J2" = 0.108006609972784
    "M2" = 7.369448
    "N2" = 14.738896
    "O2" = 0.7452608427984224
    "P2" = 0.661061693471796
    "Q2" = 186.933207938732
    "R2" = 747.732831754928
    "S2" = 0.1091671665360576
    "T2" = 0.07139903249475639
    "G3" = 25.3659715
    "H3" = 50.731943
    "I3" = 0.1464818225604603
    "J3" = 0.108006609972784
    "O3" = 0.01116592909756377
    "P3" = 0.01485661309677453
    "Q3" = 2.800741466553333
    "R3" = 16.80444879932
    "S3" = 0.001635605644792017
    "T3" = 0.001604612416259882
    "G4" = 25.3659715
    "H4" = 50.731943
    "I4" = 0.1464818225604603
    "J4" = 0.108006609972784
    "M4" = 0.084843
    "N4" = 0.254529
    "O4" = 0.008580040959044227
    "P4" = 0.0114160091622658
    "Q4" = 2.1521251199745
    "R4" = 12.912750719847
    "S4" = 0.001256820037324198
    "T4" = 0.001233004449034571
    "G5" = 25.3659715
    "H5" = 50.731943
    "I5" = 0.1464818225604603
    "J5" = 0.108006609972784
    "M5" = 2.32371
    "N5" = 6.97113
    "O5" = 0.2349931871449696
    "P5" = 0.3126656842691638
    "Q5" = 58.943161634265
    "R5" = 353.65896980559
    "S5" = 0.03442223034228648
    "T5" = 0.03376996061273321
    "I6" = 0.1917310887460603
    "J6" = 0.2120559182007618
    "M6" = 7.369448
    "N6" = 14.738896
    "O6" = 0.7452608427984224
    "P6" = 0.661061693471796
    "Q6" = 244.6781918356693
    "R6" = 1468.069151014016
    "S6" = 0.142889672789548
    "T6" = 0.1401820443965122
    "I7" = 0.1917310887460603
    "J7" = 0.2120559182007618
    "O7" = 0.01116592909756377
    "P7" = 0.01485661309677453
    "S7" = 0.002140855742737216
    "T7" = 0.003150432731589986
    "I8" = 0.1917310887460603
    "J8" = 0.2120559182007618
    "M8" = 0.084843
    "N8" = 0.254529
    "O8" = 0.008580040959044227
    "P8" = 0.0114160091622658
    "Q8" = 2.816931719976
    "R8" = 25.352385479784
    "S8" = 0.001645060594563341
    "T8" = 0.002420832305092584
    "I9" = 0.1917310887460603
    "J9" = 0.2120559182007618
    "M9" = 2.32371
    "N9" = 6.97113
    "O9" = 0.2349931871449696
    "P9" = 0.3126656842691638
    "Q9" = 77.15111920872
    "R9" = 694.3600728784801
    "S9" = 0.04505549961921174
    "T9" = 0.06630260876756702
    "G10" = 32.638883
    "H10" = 97.91664900000001
    "I10" = 0.188480976105237
    "J10" = 0.2084612710060207
    "M10" = 7.369448
    "N10" = 14.738896
    "O10" = 0.7452608427984224
    "P10" = 0.661061693471796
    "Q10" = 240.530551046584
    "R10" = 1443.183306279504
    "S10" = 0.1404674911036582
    "T10" = 0.137805760834523
    "G11" = 32.638883
    "H11" = 97.91664900000001
    "I11" = 0.188480976105237
    "J11" = 0.2084612710060207
    "O11" = 0.01116592909756377
    "P11" = 0.01485661309677453
    "Q11" = 3.603767868306666
    "R11" = 32.43391081476
    "S11" = 0.002104565215430688
    "T11" = 0.003097028448998311
    "G12" = 32.638883
    "H12" = 97.91664900000001
    "I12" = 0.188480976105237
    "J12" = 0.2084612710060207
    "M12" = 0.084843
    "N12" = 0.254529
    "O12" = 0.008580040959044227
    "P12" = 0.0114160091622658
    "Q12" = 2.769180750369
    "R12" = 24.922626753321
    "S12" = 0.00161717449498357
    "T12" = 0.002379795779782307
    "G13" = 32.638883
    "H13" = 97.91664900000001
    "I13" = 0.188480976105237
    "J13" = 0.2084612710060207
    "M13" = 2.32371
    "N13" = 6.97113
    "O13" = 0.2349931871449696
    "P13" = 0.3126656842691638
    "Q13" = 75.84329881593
    "R13" = 682.5896893433701
    "S13" = 0.04429174529116451
    "T13" = 0.06517868594271706
    "G14" = 24.426712
    "H14" = 48.853424
    "I14" = 0.1410578456622277
    "J14" = 0.1040073058467926
    "M14" = 7.369448
    "N14" = 14.738896
    "O14" = 0.7452608427984224
    "P14" = 0.661061693471796
    "Q14" = 180.011383894976
    "R14" = 720.045535579904
    "S14" = 0.1051248889415616
    "T14" = 0.06875524573651971
    "G15" = 24.426712
    "H15" = 48.853424
    "I15" = 0.1410578456622277
    "J15" = 0.1040073058467926
    "O15" = 0.01116592909756377
    "P15" = 0.01485661309677453
    "Q15" = 2.697034694293333
    "R15" = 16.18220816576
    "S15" = 0.001575041903319528
    "T15" = 0.001545196302203692
    "G16" = 24.426712
    "H16" = 48.853424
    "I16" = 0.1410578456622277
    "J16" = 0.1040073058467926
    "M16" = 0.084843
    "N16" = 0.254529
    "O16" = 0.008580040959044227
    "P16" = 0.0114160091622658
    "Q16" = 2.072435526216
    "R16" = 12.434613157296
    "S16" = 0.001210282093376452
    "T16" = 0.001187348356489565
    "G17" = 24.426712
    "H17" = 48.853424
    "I17" = 0.1410578456622277
    "J17" = 0.1040073058467926
    "M17" = 2.32371
    "N17" = 6.97113
    "O17" = 0.2349931871449696
    "P17" = 0.3126656842691638
    "Q17" = 56.76059494152001
    "R17" = 340.56356964912
    "S17" = 0.03314763272397011
    "T17" = 0.0325195154515796
    "G18" = 28.45836133333333
    "H18" = 85.37508399999999
    "I18" = 0.164339561573299
    "J18" = 0.1817606985599127
    "M18" = 7.369448
    "N18" = 14.738896
    "O18" = 0.7452608427984224
    "P18" = 0.661061693471796
    "Q18" = 209.7224140112106
    "R18" = 1258.334484067264
    "S18" = 0.12247584016324
    "T18" = 0.1201550351966325
    "G19" = 28.45836133333333
    "H19" = 85.37508399999999
    "I19" = 0.164339561573299
    "J19" = 0.1817606985599127
    "O19" = 0.01116592909756377
    "P19" = 0.01485661309677453
    "Q19" = 3.142182536017777
    "R19" = 28.27964282415999
    "S19" = 0.001835003892452172
    "T19" = 0.002700348374704086
    "G20" = 28.45836133333333
    "H20" = 85.37508399999999
    "I20" = 0.164339561573299
    "J20" = 0.1817606985599127
    "M20" = 0.084843
    "N20" = 0.254529
    "O20" = 0.008580040959044227
    "P20" = 0.0114160091622658
    "Q20" = 2.414492750603999
    "R20" = 21.730434755436
    "S20" = 0.001410040169490276
    "T20" = 0.002074981800099796
    "G21" = 28.45836133333333
    "H21" = 85.37508399999999
    "I21" = 0.164339561573299
    "J21" = 0.1817606985599127
    "M21" = 2.32371
    "N21" = 6.97113
    "O21" = 0.2349931871449696
    "P21" = 0.3126656842691638
    "Q21" = 66.12897881388
    "R21" = 595.1608093249199
    "S21" = 0.0386186773481165
    "T21" = 0.05683033318847632
    "G22" = 29.076423
    "H22" = 87.229269
    "I22" = 0.1679087053527158
    "J22" = 0.1857081964137282
    "M22" = 7.369448
    "N22" = 14.738896
    "O22" = 0.7452608427984224
    "P22" = 0.661061693471796
    "Q22" = 214.277187324504
    "R22" = 1285.663123947024
    "S22" = 0.125135783264357
    "T22" = 0.1227645748128521
    "G23" = 29.076423
    "H23" = 87.229269
    "I23" = 0.1679087053527158
    "J23" = 0.1857081964137282
    "O23" = 0.01116592909756377
    "P23" = 0.01485661309677453
    "Q23" = 3.21042478484
    "R23" = 28.89382306356
    "S23" = 0.001874856698832151
    "T23" = 0.00275899482301857
    "G24" = 29.076423
    "H24" = 87.229269
    "I24" = 0.1679087053527158
    "J24" = 0.1857081964137282
    "M24" = 0.084843
    "N24" = 0.254529
    "O24" = 0.008580040959044227
    "P24" = 0.0114160091622658
    "Q24" = 2.466930956589
    "R24" = 22.202378609301
    "S24" = 0.00144066356930639
    "T24" = 0.002120046471766977
    "G25" = 29.076423
    "H25" = 87.229269
    "I25" = 0.1679087053527158
    "J25" = 0.1857081964137282
    "M25" = 2.32371
    "N25" = 6.97113
    "O25" = 0.2349931871449696
    "P25" = 0.3126656842691638
    "Q25" = 67.56517488933001
    "R25" = 608.0865740039701
    "S25" = 0.0394574018202203
    "T25" = 0.05806458030609059
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
